$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 18
$ws.Range("H2").Value = 20

# Row 4
$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 12
$ws.Range("H4").Value = 15

# Row 5
$ws.Range("E5").Value = 34
$ws.Range("F5").Value = 17
$ws.Range("H5").Value = 21

# Row 6
$ws.Range("F6").Value = 30
$ws.Range("H6").Value = 37

# Row 7
$ws.Range("F7").Value = 22
$ws.Range("H7").Value = 23

# Row 8
$ws.Range("F8").Value = 20
$ws.Range("H8").Value = 26

# Row 10
$ws.Range("F10").Value = 17
$ws.Range("H10").Value = 19

# Row 14
$ws.Range("F14").Value = 20
$ws.Range("H14").Value = 22

# Row 15
$ws.Range("E15").Value = 117
$ws.Range("F15").Value = 58
$ws.Range("H15").Value = 69

# Row 16
$ws.Range("E16").Value = 335
$ws.Range("F16").Value = 103
$ws.Range("H16").Value = 191

# Row 17
$ws.Range("E17").Value = 32
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = 16

$wb.Save()
